# feat: add 2022-Q1 data
#
# The workbook has per-quarter "fund holders" sheets (2020-Q4 .. 2021-Q4)
# followed by a "总计" (totals) summary sheet. This adds a new "2022-Q1"
# fund-holders sheet (inserted right before "总计") and prepends a
# corresponding "2022-Q1" row to the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Remove the existing "总计" sheet - we'll recreate it at the end so
#    its tab stays last (after the newly inserted "2022-Q1" sheet).
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Delete()

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q1" fund-detail sheet right after "2021-Q4".
# ---------------------------------------------------------------------
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $prevQuarter)
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$fundRows = @(
    @("040008", "华安策略优选混合",         "52.77", "92.07", "4.63", "2.4433", 7),
    @("004495", "博时量化平衡混合",         "10.49", "38.32", "1.11", "0.1164", 8),
    @("519097", "新华中小市值优选混合",     "0.75",  "62.70", "5.17", "0.0388", 2),
    @("005616", "东方量化成长灵活配置混合", "0.20",  "92.24", "2.14", "0.0043", 6),
    @("005443", "国金量化多策略灵活配置混合","0.51", "64.10", "0.67", "0.0034", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Range("A$r").Value = ($r - 2)

    $q1.Range("B$r").NumberFormat = "@"
    $q1.Range("B$r").Value = $row[0]

    $q1.Range("C$r").Value = $row[1]

    $q1.Range("D$r").NumberFormat = "@"
    $q1.Range("D$r").Value = $row[2]

    $q1.Range("E$r").NumberFormat = "@"
    $q1.Range("E$r").Value = $row[3]

    $q1.Range("F$r").NumberFormat = "@"
    $q1.Range("F$r").Value = $row[4]

    $q1.Range("G$r").NumberFormat = "@"
    $q1.Range("G$r").Value = $row[5]

    $q1.Range("H$r").Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Recreate the "总计" sheet after "2022-Q1" with the new summary row
#    prepended in front of the pre-existing quarters.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Add($null, $q1)
$totals.Name = "总计"

$totals.Range("B1").Value = "日期"
$totals.Range("C1").Value = "持有数量(只)"
$totals.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 5, 2.61),
    @("2021-Q4", 5, 1.83),
    @("2021-Q3", 24, 9.529999999999999),
    @("2021-Q2", 11, 1.15),
    @("2021-Q1", 13, 0.75),
    @("2020-Q4", 3, 0.11)
)

$r = 2
foreach ($row in $summaryRows) {
    $totals.Range("A$r").Value = ($r - 2)
    $totals.Range("B$r").Value = $row[0]
    $totals.Range("C$r").Value = $row[1]
    $totals.Range("D$r").Value = $row[2]
    $r = $r + 1
}
